$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1886120996441281
$ws.Range("C2").Value = 0.5444839857651246
$ws.Range("J2").Value = 0.02491103202846975
$ws.Range("P2").Value = 0.1209964412811388
$ws.Range("S2").Value = 0.1209964412811388
$ws.Range("B3").Value = 0.006535947712418301
$ws.Range("C3").Value = 0.006535947712418301
$ws.Range("J3").Value = 0.0261437908496732
$ws.Range("P3").Value = 0.7320261437908496
$ws.Range("S3").Value = 0.2287581699346405
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.7368421052631579
$ws.Range("S4").Value = 0.2105263157894737
$ws.Range("B6").Value = 0.04504504504504504
$ws.Range("D6").Value = 0.01351351351351351
$ws.Range("F6").Value = 0.04504504504504504
$ws.Range("J6").Value = 0.2702702702702703
$ws.Range("O6").Value = 0.04054054054054054
$ws.Range("Q6").Value = 0.2117117117117117
$ws.Range("R6").Value = 0.05405405405405406
$ws.Range("S6").Value = 0.3198198198198198
$ws.Range("B7").Value = 0.1339712918660287
$ws.Range("D7").Value = 0.01435406698564593
$ws.Range("F7").Value = 0.03349282296650718
$ws.Range("J7").Value = 0.1148325358851675
$ws.Range("O7").Value = 0.01913875598086124
$ws.Range("Q7").Value = 0.1818181818181818
$ws.Range("R7").Value = 0.07177033492822966
$ws.Range("S7").Value = 0.430622009569378
$ws.Range("B8").Value = 0.08528784648187633
$ws.Range("D8").Value = 0.0255863539445629
$ws.Range("F8").Value = 0.0511727078891258
$ws.Range("J8").Value = 0.1407249466950959
$ws.Range("O8").Value = 0.02985074626865672
$ws.Range("Q8").Value = 0.1876332622601279
$ws.Range("R8").Value = 0.07462686567164178
$ws.Range("S8").Value = 0.4051172707889126
$ws.Range("B9").Value = 0.08
$ws.Range("D9").Value = 0.02666666666666667
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.08
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.1111111111111111
$ws.Range("S9").Value = 0.4222222222222222
$ws.Range("B10").Value = 0.112488928255093
$ws.Range("D10").Value = 0.01682905225863596
$ws.Range("E10").Value = 0.0008857395925597874
$ws.Range("F10").Value = 0.09477413640389726
$ws.Range("J10").Value = 0.1133746678476528
$ws.Range("O10").Value = 0.01771479185119575
$ws.Range("Q10").Value = 0.1886625332152347
$ws.Range("R10").Value = 0.09388839681133747
$ws.Range("S10").Value = 0.3613817537643932
$ws.Range("G11").Value = 0.1355932203389831
$ws.Range("J11").Value = 0.08135593220338982
$ws.Range("K11").Value = 0.1661016949152542
$ws.Range("L11").Value = 0.6067796610169491
$ws.Range("S11").Value = 0.01016949152542373
$ws.Range("G12").Value = 0.7604166666666666
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("K12").Value = 0.015625
$ws.Range("L12").Value = 0.046875
$ws.Range("S12").Value = 0.01041666666666667
$ws.Range("G13").Value = 0.6956521739130435
$ws.Range("J13").Value = 0.2391304347826087
$ws.Range("S13").Value = 0.06521739130434782
$ws.Range("F15").Value = 0.04090909090909091
$ws.Range("H15").Value = 0.1727272727272727
$ws.Range("I15").Value = 0.07272727272727272
$ws.Range("J15").Value = 0.2909090909090909
$ws.Range("K15").Value = 0.07272727272727272
$ws.Range("M15").Value = 0.01818181818181818
$ws.Range("O15").Value = 0.04090909090909091
$ws.Range("S15").Value = 0.2909090909090909
$ws.Range("F16").Value = 0.01169590643274854
$ws.Range("H16").Value = 0.2105263157894737
$ws.Range("I16").Value = 0.1052631578947368
$ws.Range("J16").Value = 0.2982456140350877
$ws.Range("K16").Value = 0.1052631578947368
$ws.Range("M16").Value = 0.02339181286549707
$ws.Range("O16").Value = 0.07017543859649122
$ws.Range("S16").Value = 0.1754385964912281
$ws.Range("F17").Value = 0.02331002331002331
$ws.Range("H17").Value = 0.2027972027972028
$ws.Range("I17").Value = 0.1188811188811189
$ws.Range("J17").Value = 0.3613053613053613
$ws.Range("K17").Value = 0.1188811188811189
$ws.Range("M17").Value = 0.01864801864801865
$ws.Range("O17").Value = 0.04662004662004662
$ws.Range("S17").Value = 0.1095571095571096
$ws.Range("F18").Value = 0.01036269430051814
$ws.Range("H18").Value = 0.1658031088082902
$ws.Range("I18").Value = 0.1139896373056995
$ws.Range("J18").Value = 0.3937823834196891
$ws.Range("K18").Value = 0.09326424870466321
$ws.Range("O18").Value = 0.07772020725388601
$ws.Range("S18").Value = 0.1450777202072539
$ws.Range("F19").Value = 0.01196172248803828
$ws.Range("H19").Value = 0.2256778309409888
$ws.Range("I19").Value = 0.09569377990430622
$ws.Range("J19").Value = 0.3373205741626794
$ws.Range("K19").Value = 0.1068580542264753
$ws.Range("M19").Value = 0.02551834130781499
$ws.Range("N19").Value = 0.001594896331738437
$ws.Range("O19").Value = 0.06937799043062201
